# Lower Feather recapture metadata - update to 3-6-2023 raw file revision.
#
# The authoring session removed six attribute rows that are no longer part
# of the "recapture" metadata table (lifeStage, mort, forkLength,
# totalLength, actualCountID, markCode), revised several attribute
# definitions to reflect the current data dictionary, and left the
# selection on the last row that was deleted (row 10, "actualCount").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the obsolete attribute rows. Deleting bottom-to-top keeps the
# remaining row numbers stable for each subsequent delete.
$ws.Rows.Item(24).Delete()   # markCode
$ws.Rows.Item(16).Delete()   # actualCountID
$ws.Rows.Item(14).Delete()   # totalLength
$ws.Rows.Item(13).Delete()   # forkLength
$ws.Rows.Item(11).Delete()   # mort
$ws.Rows.Item(10).Delete()   # lifeStage

# Revise attribute definitions for rows that survived the cleanup.
$ws.Range("B8").Value = 'Run revised after field visit. This is the field used in analysis. Levels = c("Not recorded", "Fall")'
$ws.Range("B16").Value = 'type of mark on fish. Levels = c("Pigment / dye", "Elastomer")'
$ws.Range("B17").Value = 'color of mark on fish. Levels = c("Brown", "Red")'
$ws.Range("B18").Value = 'position of mark on body of fish. Levels = c("whole body", "nose")'

# Leave the same selection state the editor ended the session in.
$ws.Rows.Item(10).EntireRow.Select()
